$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.344.58'
$ws.Range('E2').Value = '  +1.18%  '
$ws.Range('D3').Value = '1.858.25'
$ws.Range('E3').Value = '  +1.78%  '
$ws.Range('E4').Value = '  -0.76%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '314.27'
$c.ClearFormats()
$ws.Range('E6').Value = '  -0.66%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.4616'
$c.ClearFormats()
$ws.Range('E7').Value = '  -0.20%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.3704'
$c.ClearFormats()
$ws.Range('E8').Value = '  +0.30%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.07326'
$c.ClearFormats()
$ws.Range('E9').Value = '  +1.00%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.8828'
$c.ClearFormats()
$ws.Range('E10').Value = '  +2.46%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.07860'
$c.ClearFormats()
$ws.Range('E11').Value = '  +0.60%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '19.85'
$c.ClearFormats()
$ws.Range('E12').Value = '  -0.33%  '
$ws.Range('D13').Value = '1.878.27'
$ws.Range('E13').Value = '  +1.60%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '5.382'
$c.ClearFormats()
$ws.Range('E14').Value = '  +0.79%  '
$ws.Range('E15').Value = '  +0.41%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '91.90'
$c.ClearFormats()
$ws.Range('E16').Value = '  +0.11%  '
$ws.Range('E17').Value = '  -0.65%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '0.000008856'
$c.ClearFormats()
$ws.Range('E18').Value = '  +1.75%  '
$ws.Range('E19').Value = '  -0.65%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '14.84'
$c.ClearFormats()
$ws.Range('E20').Value = '  +2.17%  '
$ws.Range('D21').Value = '27.362.11'
$ws.Range('E21').Value = '  +0.79%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '5.117'
$c.ClearFormats()
$ws.Range('E22').Value = '  -0.64%  '
$ws.Range('E23').Value = '  -0.34%  '
$ws.Range('D24').Value = '2.128.47'
$ws.Range('E24').Value = '  +2.24%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '1.884'
$c.ClearFormats()
$ws.Range('E25').Value = '  +2.15%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '152.14'
$c.ClearFormats()
$ws.Range('E26').Value = '  -0.48%  '
$ws.Range('E27').Value = '  +0.98%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '2.079'
$c.ClearFormats()
$ws.Range('E28').Value = '  -0.64%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '5.133'
$c.ClearFormats()
$ws.Range('E29').Value = '  +0.27%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '115.99'
$c.ClearFormats()
$ws.Range('E30').Value = '  +0.54%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '0.08854'
$c.ClearFormats()
$ws.Range('E31').Value = '  +0.16%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '0.7627'
$c.ClearFormats()
$ws.Range('E32').Value = '  +5.45%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '3.021'
$c.ClearFormats()
$ws.Range('E33').Value = '  +1.83%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '1.174'
$c.ClearFormats()
$ws.Range('E34').Value = '  +3.66%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '4.490'
$c.ClearFormats()
$ws.Range('E35').Value = '  +1.09%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '2.621'
$c.ClearFormats()
$ws.Range('E36').Value = '  +6.96%  '
$ws.Range('E37').Value = '  +1.05%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '1.078'
$c.ClearFormats()
$ws.Range('E38').Value = '  -0.04%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '2.990'
$c.ClearFormats()
$ws.Range('E39').Value = '  +1.27%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '0.05203'
$c.ClearFormats()
$ws.Range('E40').Value = '  -0.66%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '7.043'
$c.ClearFormats()
$ws.Range('E41').Value = '  -2.57%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '0.5161'
$c.ClearFormats()
$ws.Range('E42').Value = '  -0.05%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '0.1640'
$c.ClearFormats()
$ws.Range('E43').Value = '  +0.77%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '8.350'
$c.ClearFormats()
$ws.Range('E44').Value = '  +2.00%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '0.4835'
$c.ClearFormats()
$ws.Range('E45').Value = '  +0.58%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '10.32'
$c.ClearFormats()
$ws.Range('E46').Value = '  +1.32%  '
$ws.Range('E47').Value = '  -0.73%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '103.36'
$c.ClearFormats()
$ws.Range('E48').Value = '  +0.60%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '1.652'
$c.ClearFormats()
$ws.Range('E49').Value = '  +2.16%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '0.06220'
$c.ClearFormats()
$ws.Range('E50').Value = '  -0.70%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '65.67'
$c.ClearFormats()
$ws.Range('E51').Value = '  +2.01%  '
